$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @("Qwen2-5B-Instruct-simple 12452", "0.2 ± 0.04", "0.98 ± 0.08", "0.06 ± 0.49", "0.21 ± 0.58", "0.4 ± 0.13", "0.003 ± 0.00"),
    @("Qwen2-5B-Instruct-perso old", "0.46 ± 0.29", "0.92 ± 0.18", "-0.01 ± 0.46", "0.27 ± 0.65", "0.31 ± 0.12", "0.0 ± 0.00"),
    @("Qwen2-5B-Instruct-perso check-2304", "0.51 ± 0.25", "0.71 ± 0.44", "-0.13 ± 0.62", "0.17 ± 0.52", "0.32 ± 0.2", "0.221 ± 0.00"),
    @("Qwen2-5B-Instruct-perso", "0.89 ± 0.04", "0.7 ± 0.45", "0.18 ± 0.87", "0.3 ± 0.62", "0.44 ± 0.3", "0.29 ± 0.00")
)

$startRow = 7
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $rowValues[$j]
    }
}
